$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1880566.9
$ws.Range("I19").Value = 4386595
$ws.Range("J19").Value = 1045.875
$ws.Range("K19").Value = 4386595
$ws.Range("L19").Value = 1045.875
$ws.Range("M19").Value = -4386420
$ws.Range("N19").Value = -1395.875
$ws.Range("H40").Value = 1252.0667
$ws.Range("I40").Value = 974.875
$ws.Range("J40").Value = 1568.8572
$ws.Range("K40").Value = 974.875
$ws.Range("L40").Value = 1568.8572
$ws.Range("M40").Value = -799.875
$ws.Range("N40").Value = -1918.8572
$ws.Range("H43").Value = 2084.2222
$ws.Range("I43").Value = 1323.3334
$ws.Range("K43").Value = 1323.3334
$ws.Range("M43").Value = -1254.3334
$ws.Range("H115").Value = 1415.909
$ws.Range("J115").Value = 3000
$ws.Range("L115").Value = 9000
$ws.Range("N115").Value = -12134
$ws.Range("H129").Value = 855.08
$ws.Range("J129").Value = 875.11456
$ws.Range("L129").Value = 2625.34368
$ws.Range("N129").Value = -12625.34368
$ws.Range("H141").Value = 15116.134
$ws.Range("I141").Value = 16672.46
$ws.Range("K141").Value = 50017.38
$ws.Range("M141").Value = -44837.38

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 649
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 649
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 649
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -875
$ws.Range("H24").Value = 24450.7
$ws.Range("J24").Value = 24450.7
$ws.Range("L24").Value = 24450.7
$ws.Range("N24").Value = -25198.7
$ws.Range("H32").Value = 5596.8477
$ws.Range("I32").Value = 4561.375
$ws.Range("J32").Value = 12500
$ws.Range("K32").Value = 4561.375
$ws.Range("L32").Value = 12500
$ws.Range("M32").Value = -4274.375
$ws.Range("N32").Value = -13074
$ws.Range("H75").Value = 39800
$ws.Range("J75").Value = 39800
$ws.Range("L75").Value = 39800
$ws.Range("N75").Value = -41548
$ws.Range("H78").Value = 39800
$ws.Range("J78").Value = 39800
$ws.Range("L78").Value = 119400
$ws.Range("N78").Value = -128136
$ws.Range("H100").Value = 24450.7
$ws.Range("J100").Value = 24450.7
$ws.Range("L100").Value = 24450.7
$ws.Range("N100").Value = -26614.7
$ws.Range("H116").Value = 649
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 649
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 649
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -5237
$ws.Range("H122").Value = 3515.5789
$ws.Range("I122").Value = 1372.7693
$ws.Range("K122").Value = 4118.3079
$ws.Range("M122").Value = -1668.3079
$ws.Range("H137").Value = 38337.5
$ws.Range("I137").Value = 29800
$ws.Range("J137").Value = 41183.332
$ws.Range("K137").Value = 29800
$ws.Range("L137").Value = 41183.332
$ws.Range("M137").Value = -24700
$ws.Range("N137").Value = -51383.332

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 649
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 649
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 649
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -877
$ws.Range("H59").Value = 43333.332
$ws.Range("J59").Value = 43333.332
$ws.Range("L59").Value = 43333.332
$ws.Range("N59").Value = -45027.332
$ws.Range("H137").Value = 45700
$ws.Range("J137").Value = 45700
$ws.Range("L137").Value = 45700
$ws.Range("N137").Value = -55900

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3129.746
$ws.Range("I58").Value = 1807.3334
$ws.Range("K58").Value = 1807.3334
$ws.Range("M58").Value = -1604.3334
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524
$ws.Range("H132").Value = 5250
$ws.Range("I132").Value = 3306
$ws.Range("K132").Value = 9918
$ws.Range("M132").Value = -7388
$ws.Range("H134").Value = 4142.7144
$ws.Range("I134").Value = 2000
$ws.Range("K134").Value = 6000
$ws.Range("M134").Value = -3465
$ws.Range("H136").Value = 3129.746
$ws.Range("I136").Value = 1807.3334
$ws.Range("K136").Value = 5422.0002
$ws.Range("M136").Value = -2872.0002
$ws.Range("H137").Value = 45238.57
$ws.Range("J137").Value = 45238.57
$ws.Range("L137").Value = 45238.57
$ws.Range("N137").Value = -55438.57

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 17171.428
$ws.Range("J39").Value = 18700
$ws.Range("L39").Value = 56100
$ws.Range("N39").Value = -56688
$ws.Range("H68").Value = 1268.3334
$ws.Range("I68").Value = 602
$ws.Range("J68").Value = 1601.5
$ws.Range("K68").Value = 1806
$ws.Range("L68").Value = 4804.5
$ws.Range("M68").Value = -995
$ws.Range("N68").Value = -6426.5
$ws.Range("H71").Value = 1268.3334
$ws.Range("I71").Value = 602
$ws.Range("J71").Value = 1601.5
$ws.Range("K71").Value = 5418
$ws.Range("L71").Value = 14413.5
$ws.Range("M71").Value = -1362
$ws.Range("N71").Value = -22525.5
$ws.Range("H131").Value = 699.08
$ws.Range("I131").Value = 273.90475
$ws.Range("J131").Value = 812.10126
$ws.Range("K131").Value = 821.71425
$ws.Range("L131").Value = 2436.30378
$ws.Range("M131").Value = 4218.28575
$ws.Range("N131").Value = -12516.30378

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 35151.2
$ws.Range("J46").Value = 35151.2
$ws.Range("L46").Value = 35151.2
$ws.Range("N46").Value = -35463.2
$ws.Range("H126").Value = 3342.88
$ws.Range("I126").Value = 2880.3948
$ws.Range("J126").Value = 4807.4165
$ws.Range("K126").Value = 8641.1844
$ws.Range("L126").Value = 14422.2495
$ws.Range("M126").Value = -6171.1844
$ws.Range("N126").Value = -19362.2495
$ws.Range("H130").Value = 47225.9
$ws.Range("J130").Value = 47225.9
$ws.Range("L130").Value = 47225.9
$ws.Range("N130").Value = -57265.9
$ws.Range("H137").Value = 40780
$ws.Range("J137").Value = 40780
$ws.Range("L137").Value = 40780
$ws.Range("N137").Value = -50980

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2386.25
$ws.Range("I61").Value = 2619
$ws.Range("J61").Value = 2060.4
$ws.Range("K61").Value = 2619
$ws.Range("L61").Value = 2060.4
$ws.Range("M61").Value = -2417
$ws.Range("N61").Value = -2464.4
$ws.Range("H92").Value = 32161.5
$ws.Range("J92").Value = 32161.5
$ws.Range("L92").Value = 32161.5
$ws.Range("N92").Value = -37153.5
$ws.Range("H113").Value = 2386.25
$ws.Range("I113").Value = 2619
$ws.Range("J113").Value = 2060.4
$ws.Range("K113").Value = 2619
$ws.Range("L113").Value = 2060.4
$ws.Range("M113").Value = -449
$ws.Range("N113").Value = -6400.4
$ws.Range("H132").Value = 5731.5835
$ws.Range("I132").Value = 3379.6667
$ws.Range("J132").Value = 8083.5
$ws.Range("K132").Value = 10139.0001
$ws.Range("L132").Value = 24250.5
$ws.Range("M132").Value = -7609.000100000001
$ws.Range("N132").Value = -29310.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3255.0435
$ws.Range("I122").Value = 1414.8572
$ws.Range("J122").Value = 6117.5557
$ws.Range("K122").Value = 4244.571599999999
$ws.Range("L122").Value = 18352.6671
$ws.Range("M122").Value = -1794.571599999999
$ws.Range("N122").Value = -23252.6671
$ws.Range("H136").Value = 9219
$ws.Range("J136").Value = 11753.286
$ws.Range("L136").Value = 35259.858
$ws.Range("N136").Value = -40359.858
